$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 293.66666
$ws.Range("I2").Value = 309
$ws.Range("K2").Value = 309
$ws.Range("M2").Value = -196
$ws.Range("H17").Value = 1069.6586
$ws.Range("J17").Value = 1069.6586
$ws.Range("L17").Value = 3208.9758
$ws.Range("N17").Value = -3544.9758
$ws.Range("H46").Value = 34512.5
$ws.Range("I46").Value = 45100
$ws.Range("J46").Value = 2750
$ws.Range("K46").Value = 135300
$ws.Range("L46").Value = 8250
$ws.Range("M46").Value = -135181
$ws.Range("N46").Value = -8488
$ws.Range("H59").Value = 3000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H60").Value = 34512.5
$ws.Range("I60").Value = 45100
$ws.Range("J60").Value = 2750
$ws.Range("K60").Value = 135300
$ws.Range("L60").Value = 8250
$ws.Range("M60").Value = -134816
$ws.Range("N60").Value = -9218
$ws.Range("H61").Value = 750
$ws.Range("I61").Value = 750
$ws.Range("K61").Value = 2250
$ws.Range("M61").Value = -2078
$ws.Range("H98").Value = 7825.885
$ws.Range("I98").Value = 8238.087
$ws.Range("K98").Value = 8238.087
$ws.Range("M98").Value = -6740.087
$ws.Range("H115").Value = 3887.1333
$ws.Range("I115").Value = 686.7143
$ws.Range("K115").Value = 2060.1429
$ws.Range("M115").Value = -493.1428999999998
$ws.Range("H122").Value = 7825.885
$ws.Range("I122").Value = 8238.087
$ws.Range("K122").Value = 24714.261
$ws.Range("M122").Value = -22264.261
$ws.Range("H138").Value = 2110.1099
$ws.Range("I138").Value = 1069.1428
$ws.Range("K138").Value = 3207.4284
$ws.Range("M138").Value = 1932.5716
$ws.Range("H141").Value = 6493.593
$ws.Range("I141").Value = 4287.7144
$ws.Range("K141").Value = 12863.1432
$ws.Range("M141").Value = -7683.143199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3860.4546
$ws.Range("I2").Value = 823.3333
$ws.Range("J2").Value = 4999.375
$ws.Range("K2").Value = 823.3333
$ws.Range("L2").Value = 4999.375
$ws.Range("M2").Value = -710.3333
$ws.Range("N2").Value = -5225.375
$ws.Range("H21").Value = 6999.5
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10748
$ws.Range("H32").Value = 2625.4893
$ws.Range("I32").Value = 2719.9556
$ws.Range("K32").Value = 2719.9556
$ws.Range("M32").Value = -2432.9556
$ws.Range("H74").Value = 2240.45
$ws.Range("I74").Value = 1958.5555
$ws.Range("K74").Value = 1958.5555
$ws.Range("M74").Value = -1084.5555
$ws.Range("H77").Value = 2240.45
$ws.Range("I77").Value = 1958.5555
$ws.Range("K77").Value = 9792.7775
$ws.Range("M77").Value = -5424.7775
$ws.Range("H116").Value = 3860.4546
$ws.Range("I116").Value = 823.3333
$ws.Range("J116").Value = 4999.375
$ws.Range("K116").Value = 823.3333
$ws.Range("L116").Value = 4999.375
$ws.Range("M116").Value = 1470.6667
$ws.Range("N116").Value = -9587.375
$ws.Range("H132").Value = 4185.391
$ws.Range("I132").Value = 2465.6924
$ws.Range("K132").Value = 7397.0772
$ws.Range("M132").Value = -4867.0772

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3860.4546
$ws.Range("I3").Value = 823.3333
$ws.Range("J3").Value = 4999.375
$ws.Range("K3").Value = 823.3333
$ws.Range("L3").Value = 4999.375
$ws.Range("M3").Value = -709.3333
$ws.Range("N3").Value = -5227.375
$ws.Range("H5").Value = 542
$ws.Range("I5").Value = 790
$ws.Range("J5").Value = 294
$ws.Range("K5").Value = 790
$ws.Range("L5").Value = 294
$ws.Range("M5").Value = -677
$ws.Range("N5").Value = -520
$ws.Range("H26").Value = 39874.145
$ws.Range("I26").Value = 40266.5
$ws.Range("J26").Value = 37520
$ws.Range("K26").Value = 40266.5
$ws.Range("L26").Value = 37520
$ws.Range("M26").Value = -39974.5
$ws.Range("N26").Value = -38104
$ws.Range("H134").Value = 2837.2
$ws.Range("I134").Value = 2828.6316
$ws.Range("K134").Value = 8485.8948
$ws.Range("M134").Value = -5950.8948

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4710.15
$ws.Range("I31").Value = 1860.4
$ws.Range("K31").Value = 1860.4
$ws.Range("M31").Value = -1565.4
$ws.Range("H34").Value = 4710.15
$ws.Range("I34").Value = 1860.4
$ws.Range("K34").Value = 1860.4
$ws.Range("M34").Value = -1658.4
$ws.Range("H55").Value = 12998
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 12998
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 12998
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -13628
$ws.Range("H134").Value = 2210
$ws.Range("I134").Value = 2079.2
$ws.Range("K134").Value = 6237.599999999999
$ws.Range("M134").Value = -3702.599999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31.708334
$ws.Range("I2").Value = 26.35
$ws.Range("J2").Value = 58.5
$ws.Range("K2").Value = 158.1
$ws.Range("L2").Value = 351
$ws.Range("M2").Value = -45.10000000000002
$ws.Range("N2").Value = -577
$ws.Range("H51").Value = 1283.1666
$ws.Range("I51").Value = 199
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 597
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -137
$ws.Range("N51").Value = -5420
$ws.Range("H134").Value = 2506.0557
$ws.Range("I134").Value = 1807.2667
$ws.Range("K134").Value = 5421.800099999999
$ws.Range("M134").Value = -351.8000999999995
$ws.Range("H138").Value = 5293.2856
$ws.Range("I138").Value = 2289.75
$ws.Range("K138").Value = 6869.25
$ws.Range("M138").Value = -1729.25
$ws.Range("H139").Value = 4222
$ws.Range("I139").Value = 1749.25
$ws.Range("J139").Value = 4928.5
$ws.Range("K139").Value = 5247.75
$ws.Range("L139").Value = 14785.5
$ws.Range("M139").Value = -107.75
$ws.Range("N139").Value = -25065.5
$ws.Range("H141").Value = 1853.1111
$ws.Range("I141").Value = 1853.1111
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5559.3333
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -379.3333000000002
$ws.Range("N141").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 20555.111
$ws.Range("I55").Value = 16000
$ws.Range("J55").Value = 22832.666
$ws.Range("K55").Value = 16000
$ws.Range("L55").Value = 22832.666
$ws.Range("M55").Value = -15673
$ws.Range("N55").Value = -23486.666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4444.1577
$ws.Range("I40").Value = 2979.9565
$ws.Range("J40").Value = 6689.2666
$ws.Range("K40").Value = 2979.9565
$ws.Range("M40").Value = -2843.9565
$ws.Range("N40").Value = -6961.2666
$ws.Range("H59").Value = 20750
$ws.Range("J59").Value = 20750
$ws.Range("L59").Value = 20750
$ws.Range("N59").Value = -22058
$ws.Range("H136").Value = 3733.1904
$ws.Range("I136").Value = 3774.85
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 11324.55
$ws.Range("L136").Value = 8700
$ws.Range("M136").Value = -8774.549999999999
$ws.Range("N136").Value = -13800

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 16663.334
$ws.Range("I24").Value = 16663.334
$ws.Range("K24").Value = 16663.334
$ws.Range("M24").Value = -16433.334
$ws.Range("H62").Value = 11050.167
$ws.Range("I62").Value = 10120.4
$ws.Range("K62").Value = 10120.4
$ws.Range("M62").Value = -9496.4
$ws.Range("H65").Value = 11050.167
$ws.Range("I65").Value = 10120.4
$ws.Range("K65").Value = 50602
$ws.Range("M65").Value = -47482
$ws.Range("H126").Value = 1795.8
$ws.Range("I126").Value = 1744.75
$ws.Range("K126").Value = 5234.25
$ws.Range("M126").Value = -2764.25
